$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.625.65'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.906.40'
$ws.Range('E3').Value = '  -2.09%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.74'
$ws.Range('E5').Value = '  -1.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.19'
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.505'
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.906.56'
$ws.Range('E9').Value = '  -1.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.93'
$ws.Range('E10').Value = '  -5.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.149'
$ws.Range('E11').Value = '  +3.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.434'
$ws.Range('E12').Value = '  -3.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000235'
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.78'
$ws.Range('E14').Value = '  -2.03%  '
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.388.32'
$ws.Range('E16').Value = '  -2.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.657.45'
$ws.Range('E17').Value = '  -1.66%  '
$ws.Range('E18').Value = '  -1.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.908.68'
$ws.Range('E19').Value = '  -2.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '436.71'
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('E21').Value = '  -0.93%  '
$ws.Range('E22').Value = '  -2.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.93'
$ws.Range('E23').Value = '  -2.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.80'
$ws.Range('E24').Value = '  -1.31%  '
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.11'
$ws.Range('E26').Value = '  -8.83%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').Value = '  -3.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000104'
$ws.Range('E29').Value = '  +18.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.18'
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.55'
$ws.Range('E31').Value = '  -2.47%  '
$ws.Range('E32').Value = '  -1.23%  '
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '25.81'
$ws.Range('E35').Value = '  -2.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.973'
$ws.Range('E36').Value = '  -1.88%  '
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.49'
$ws.Range('E37').Value = '  -2.93%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.02'
$ws.Range('E38').Value = '  +3.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.04'
$ws.Range('E39').Value = '  -1.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.99'
$ws.Range('E40').Value = '  -2.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.34'
$ws.Range('E41').Value = '  -3.25%  '
$ws.Range('E42').Value = '  -2.07%  '
$ws.Range('E43').Value = '  -4.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.71'
$ws.Range('E44').Value = '  -0.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.694.57'
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '133.67'
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '340.93'
$ws.Range('E49').Value = '  -7.72%  '
$ws.Range('E50').Value = '  -1.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.26'
$ws.Range('E51').Value = '  -4.24%  '
